# Apply the "trend summary" refresh:
#  - Summary Table: drop the placeholder dashed row, tidy header text/style,
#    restore default page margins.
#  - Cooccurrence: populate source/target/count keyword co-occurrence data.
#  - Associations: replace the placeholder row with the full term/count list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary Table
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary Table")

# Remove the "---------" placeholder row (row 2); real data shifts up.
$summary.Rows.Item(2).Delete()

# Re-write the header labels without the leading/trailing padding spaces,
# and copy the bold/bordered "s=2" header style from a sheet that already
# uses it (Associations!A1) onto the Summary Table header row.
$summary.Range("A1").Value = "Keyword"
$summary.Range("B1").Value = "Keyword Count"
$summary.Range("C1").Value = "Short Summary"
$summary.Range("D1").Value = "Source URL"
$summary.Range("E1").Value = "Detailed Summary"

$assoc = $wb.Worksheets.Item("Associations")
$assoc.Range("A1").Copy()
$summary.Range("A1:E1").PasteSpecial(-4122)

# Restore the default page margins (in points: 0.75/0.75/1/1/0.5/0.5 in).
$summary.PageSetup.LeftMargin = 54
$summary.PageSetup.RightMargin = 54
$summary.PageSetup.TopMargin = 72
$summary.PageSetup.BottomMargin = 72
$summary.PageSetup.HeaderMargin = 36
$summary.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Cooccurrence
# ---------------------------------------------------------------------------
$cooc = $wb.Worksheets.Item("Cooccurrence")

$cooc.Range("A1").Value = "source"
$cooc.Range("B1").Value = "target"
$cooc.Range("C1").Value = "count"

$coocData = @(
    @("人工智能", "新质生产力", 2),
    @("人工智能", "生物技术", 2),
    @("人工智能", "科技成果转化", 2),
    @("新质生产力", "生物技术", 1),
    @("新质生产力", "科技成果转化", 1),
    @("生物技术", "科技成果转化", 1),
    @("创新驱动发展", "国家创新体系", 1),
    @("技术转移", "科技成果转化", 1)
)

$r = 2
foreach ($row in $coocData) {
    $cooc.Cells.Item($r, 1).Value = $row[0]
    $cooc.Cells.Item($r, 2).Value = $row[1]
    $cooc.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$assoc.Range("A1:B1").Copy()
$cooc.Range("A1:C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Associations
# ---------------------------------------------------------------------------
$assocData = @(
    @("新质生产力", 2),
    @("人工智能", 7),
    @("科技成果转化", 3),
    @("生物技术", 4),
    @("量子计算", 1),
    @("氢能", 2),
    @("量子通信", 1),
    @("创新驱动发展", 1),
    @("国家创新体系", 2),
    @("科技企业孵化", 1),
    @("技术转移", 2),
    @("知识产权保护", 1)
)

$r = 2
foreach ($row in $assocData) {
    $assoc.Cells.Item($r, 1).Value = $row[0]
    $assoc.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
